$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.529.32"
$ws.Range("E2").Value = "  +5.81%  "
$ws.Range("D3").Value = "1.723.95"
$ws.Range("E3").Value = "  +4.77%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5344"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.31%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2674"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06591"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.00%  "
$ws.Range("E10").Value = "  +7.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07708"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "1.726.25"
$ws.Range("E13").Value = "  +5.45%  "
$ws.Range("D14").Value = "1.961.26"
$ws.Range("E14").Value = "  +4.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5828"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.05%  "
$ws.Range("D16").Value = "0.0₅8276"
$ws.Range("E16").Value = "  +2.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.82%  "
$ws.Range("D18").Value = "27.529.96"
$ws.Range("E18").Value = "  +5.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +15.39%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.734"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.32%  "
$ws.Range("E22").Value = "  +2.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.084"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.58%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.731"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +14.62%  "
$ws.Range("E27").Value = "  +4.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.409"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05544"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.569"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.448"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.664"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.856"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9665"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.14%  "
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5981"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01652"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.912"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8559"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.12%  "
$ws.Range("D42").Value = "1.056.68"
$ws.Range("E42").Value = "  +2.81%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("D45").Value = "1.867.75"
$ws.Range("E45").Value = "  +4.73%  "
$ws.Range("E46").Value = "  +2.95%  "
$ws.Range("E47").Value = "  +3.72%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.233"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.61%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4455"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.25%  "
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05239"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.63%  "
